$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.033.88"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "'2.243.35"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'315.62"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "'99.26"
$ws.Range("E6").Value = "  -6.65%  "
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -7.21%  "
$ws.Range("D10").Value = "'36.25"
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("D11").Value = "'0.0822"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("E12").Value = "  -7.31%  "
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").Value = "'2.587.12"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("E15").Value = "  -5.25%  "
$ws.Range("D16").Value = "'2.240.36"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "'13.97"
$ws.Range("E17").Value = "  -4.90%  "
$ws.Range("D18").Value = "'43.933.40"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = "  -7.16%  "
$ws.Range("D20").Value = "'0.0₃0977"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").Value = "'6.33"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").Value = "'65.56"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'236.45"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  -7.34%  "
$ws.Range("E25").Value = "  -8.97%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'10.16"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("D29").Value = "'36.50"
$ws.Range("E29").Value = "  -5.31%  "
$ws.Range("E30").Value = "  -9.04%  "
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").Value = "'155.75"
$ws.Range("E32").Value = "  -4.64%  "
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "'1.90"
$ws.Range("E36").Value = "  -7.91%  "
$ws.Range("E37").Value = "  -8.21%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "'15.52"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'3.53"
$ws.Range("E40").Value = "  -11.66%  "
$ws.Range("D41").Value = "'4.00"
$ws.Range("E41").Value = "  -10.99%  "
$ws.Range("E42").Value = "  -6.73%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'1.700.16"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "'82.97"
$ws.Range("E45").Value = "  -4.81%  "
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").Value = "'5.17"
$ws.Range("E47").Value = "  -6.37%  "
$ws.Range("D48").Value = "'101.82"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").Value = "'71.30"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("E50").Value = "  -6.97%  "
$ws.Range("E51").Value = "  -7.00%  "
